# Automatic update by jenkins
# Renames "<first> doe" style display names to "<first>_doe" and introduces a
# new "children" (mref) attribute on patients, together with the related
# worksheet bookkeeping (selections / active sheet).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: root_hospital_cities -> "new york" becomes "new_york"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2, 1).Value = "new_york"

# ---------------------------------------------------------------------------
# Sheet 2: root_hospital_patients
#   - displayName values switch from "x doe" to "x_doe"
#   - birthplace values switch from "new york" to "new_york"
#   - a new "children" column is inserted before "disease"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Cells.Item(2, 1).Value = "john_doe"
$ws2.Cells.Item(3, 1).Value = "jane_doe"
$ws2.Cells.Item(4, 1).Value = "papa_doe"

$ws2.Cells.Item(2, 5).Value = "new_york"
$ws2.Cells.Item(4, 5).Value = "new_york"

# Insert the new "children" column in front of the existing "disease" column (F)
$ws2.Columns.Item(6).Insert()
$ws2.Cells.Item(1, 6).Value = "children"
$ws2.Cells.Item(4, 6).Value = "john_doe, jane_doe"

# ---------------------------------------------------------------------------
# Sheet 3: root_hospital_users -> displayName references updated to match
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(2, 3).Value = "john_doe"
$ws3.Cells.Item(3, 3).Value = "jane_doe"

# ---------------------------------------------------------------------------
# Sheet 4: attributes -> add the "children" mref attribute for patients
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Rows.Item(10).Insert()
$ws4.Cells.Item(10, 1).Value = "children"
$ws4.Cells.Item(10, 2).Value = "root_hospital_patients"
$ws4.Cells.Item(10, 3).Value = "mref"
$ws4.Cells.Item(10, 5).Value = "root_hospital_patients"
$ws4.Cells.Item(10, 7).Value = "children of a patient"

# ---------------------------------------------------------------------------
# Restore per-sheet cursor positions and the active sheet/tab
# ---------------------------------------------------------------------------
$win = $excel.Windows.Item(1)
$win.TabRatio = 679

$ws1.Activate()
$ws1.Range("E35").Select()

$ws2.Activate()
$ws2.Range("G15:G16").Select()

$ws3.Activate()
$ws3.Range("E2").Select()

$ws5 = $wb.Worksheets.Item(5)
$ws5.Activate()
$ws5.Range("E21").Select()

$ws6 = $wb.Worksheets.Item(6)
$ws6.Activate()
$ws6.Range("A2").Select()

$ws4.Activate()
$ws4.Range("E28").Select()
